$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 8259.177
$ws.Range("I74").Value = 5634
$ws.Range("J74").Value = 11212.5
$ws.Range("K74").Value = 5634
$ws.Range("L74").Value = 11212.5
$ws.Range("M74").Value = -4698
$ws.Range("N74").Value = -13084.5
$ws.Range("H77").Value = 8259.177
$ws.Range("I77").Value = 5634
$ws.Range("J77").Value = 11212.5
$ws.Range("K77").Value = 28170
$ws.Range("L77").Value = 56062.5
$ws.Range("M77").Value = -23490
$ws.Range("N77").Value = -65422.5
$ws.Range("H106").Value = 6129.125
$ws.Range("I106").Value = 8001.6665
$ws.Range("K106").Value = 8001.6665
$ws.Range("M106").Value = -7370.6665
$ws.Range("H137").Value = 3881.75
$ws.Range("I137").Value = 2397.7896
$ws.Range("J137").Value = 7014.5557
$ws.Range("K137").Value = 7193.3688
$ws.Range("L137").Value = 21043.6671
$ws.Range("M137").Value = -4643.3688
$ws.Range("N137").Value = -26143.6671
$ws.Range("H138").Value = 6419.4
$ws.Range("I138").Value = 5489.6665
$ws.Range("J138").Value = 6546.1816
$ws.Range("K138").Value = 16468.9995
$ws.Range("L138").Value = 19638.5448
$ws.Range("M138").Value = -11328.9995
$ws.Range("N138").Value = -29918.5448

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4546.75
$ws.Range("I32").Value = 4292.602
$ws.Range("J32").Value = 17000
$ws.Range("K32").Value = 4292.602
$ws.Range("L32").Value = 17000
$ws.Range("M32").Value = -4005.602
$ws.Range("N32").Value = -17574
$ws.Range("H61").Value = 5139.25
$ws.Range("I61").Value = 3353.2727
$ws.Range("K61").Value = 3353.2727
$ws.Range("M61").Value = -3141.2727
$ws.Range("H74").Value = 2521.9524
$ws.Range("I74").Value = 2455.8948
$ws.Range("K74").Value = 2455.8948
$ws.Range("M74").Value = -1581.8948
$ws.Range("H77").Value = 2521.9524
$ws.Range("I77").Value = 2455.8948
$ws.Range("K77").Value = 12279.474
$ws.Range("M77").Value = -7911.474
$ws.Range("H122").Value = 3112.6858
$ws.Range("I122").Value = 2297.926
$ws.Range("K122").Value = 6893.778
$ws.Range("M122").Value = -4443.778
$ws.Range("H132").Value = 4638.838
$ws.Range("I132").Value = 4110.212
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 12330.636
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -9800.636000000002
$ws.Range("N132").Value = -32060
$ws.Range("H136").Value = 5139.25
$ws.Range("I136").Value = 3353.2727
$ws.Range("K136").Value = 10059.8181
$ws.Range("M136").Value = -7509.8181

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2642
$ws.Range("I105").Value = 1744.8572
$ws.Range("K105").Value = 1744.8572
$ws.Range("M105").Value = 2.142800000000079

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 101.27273
$ws.Range("J25").Value = 342
$ws.Range("L25").Value = 342
$ws.Range("N25").Value = -690
$ws.Range("H31").Value = 3437.38
$ws.Range("I31").Value = 2563.5
$ws.Range("J31").Value = 4748.2
$ws.Range("K31").Value = 2563.5
$ws.Range("L31").Value = 4748.2
$ws.Range("M31").Value = -2268.5
$ws.Range("N31").Value = -5338.2
$ws.Range("H34").Value = 3437.38
$ws.Range("I34").Value = 2563.5
$ws.Range("J34").Value = 4748.2
$ws.Range("K34").Value = 2563.5
$ws.Range("L34").Value = 4748.2
$ws.Range("M34").Value = -2361.5
$ws.Range("N34").Value = -5152.2
$ws.Range("H141").Value = 337087.2
$ws.Range("J141").Value = 337087.2
$ws.Range("L141").Value = 337087.2
$ws.Range("N141").Value = -347447.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 948.4666999999999
$ws.Range("J60").Value = 1253.9
$ws.Range("L60").Value = 3761.7
$ws.Range("N60").Value = -4263.700000000001
$ws.Range("H131").Value = 3598.8
$ws.Range("I131").Value = 2322.1538
$ws.Range("K131").Value = 6966.4614
$ws.Range("M131").Value = -1926.4614

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 11905653
$ws.Range("I97").Value = 886.6
$ws.Range("J97").Value = 41667572
$ws.Range("K97").Value = 886.6
$ws.Range("L97").Value = 41667572
$ws.Range("M97").Value = -390.6
$ws.Range("N97").Value = -41668564
$ws.Range("H122").Value = 13137.046
$ws.Range("I122").Value = 19152.166
$ws.Range("K122").Value = 57456.49800000001
$ws.Range("M122").Value = -55006.49800000001
$ws.Range("H132").Value = 4467.109
$ws.Range("I132").Value = 4525.359
$ws.Range("K132").Value = 13576.077
$ws.Range("M132").Value = -11046.077

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4883.0713
$ws.Range("I7").Value = 2859.6667
$ws.Range("K7").Value = 2859.6667
$ws.Range("M7").Value = -2747.6667
$ws.Range("H22").Value = 1830.4286
$ws.Range("I22").Value = 1382.2
$ws.Range("J22").Value = 2951
$ws.Range("K22").Value = 1382.2
$ws.Range("L22").Value = 2951
$ws.Range("M22").Value = -1087.2
$ws.Range("N22").Value = -3541
$ws.Range("H27").Value = 1830.4286
$ws.Range("I27").Value = 1382.2
$ws.Range("J27").Value = 2951
$ws.Range("K27").Value = 1382.2
$ws.Range("L27").Value = 2951
$ws.Range("M27").Value = -1275.2
$ws.Range("N27").Value = -3165
$ws.Range("H46").Value = 7827.9062
$ws.Range("I46").Value = 1899.5
$ws.Range("J46").Value = 8674.821
$ws.Range("K46").Value = 1899.5
$ws.Range("L46").Value = 8674.821
$ws.Range("M46").Value = -1711.5
$ws.Range("N46").Value = -9050.821
$ws.Range("H126").Value = 4883.0713
$ws.Range("I126").Value = 2859.6667
$ws.Range("K126").Value = 8579.000100000001
$ws.Range("M126").Value = -6109.000100000001
$ws.Range("H132").Value = 7033.933
$ws.Range("I132").Value = 5643.5713
$ws.Range("J132").Value = 8250.5
$ws.Range("K132").Value = 16930.7139
$ws.Range("L132").Value = 24751.5
$ws.Range("M132").Value = -14400.7139
$ws.Range("N132").Value = -29811.5
$ws.Range("H136").Value = 5280.5356
$ws.Range("I136").Value = 3350.0588
$ws.Range("J136").Value = 8264
$ws.Range("K136").Value = 10050.1764
$ws.Range("L136").Value = 24792
$ws.Range("M136").Value = -7500.1764
$ws.Range("N136").Value = -29892

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 3768066
$ws.Range("I52").Value = 6012006
$ws.Range("K52").Value = 6012006
$ws.Range("M52").Value = -6011780
$ws.Range("H132").Value = 2456.7942
$ws.Range("I132").Value = 1340.3928
$ws.Range("K132").Value = 4021.1784
$ws.Range("M132").Value = -1491.1784
